# Update "想去人数" (F column) figures across sheets to match the newly
# generated site output (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value  = 6526
$ws.Range("F4").Value  = 752
$ws.Range("F5").Value  = 1097
$ws.Range("F6").Value  = 102
$ws.Range("F7").Value  = 587
$ws.Range("F8").Value  = 207
$ws.Range("F9").Value  = 33
$ws.Range("F10").Value = 762
$ws.Range("F11").Value = 1235
$ws.Range("F12").Value = 17
$ws.Range("F14").Value = 212
$ws.Range("F15").Value = 478
$ws.Range("F16").Value = 352
$ws.Range("F17").Value = 319
$ws.Range("F18").Value = 1436
$ws.Range("F19").Value = 692
$ws.Range("F20").Value = 412
$ws.Range("F21").Value = 416
$ws.Range("F24").Value = 184
$ws.Range("F25").Value = 2258
$ws.Range("F27").Value = 133
$ws.Range("F28").Value = 413
$ws.Range("F30").Value = 3679
$ws.Range("F31").Value = 50
$ws.Range("F32").Value = 672

# --- Sheet: 演出 (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F10").Value = 1026
$ws.Range("F12").Value = 124
$ws.Range("F16").Value = 83
$ws.Range("F20").Value = 4102
$ws.Range("F24").Value = 36
$ws.Range("F25").Value = 209
$ws.Range("F35").Value = 7

# --- Sheet: 本地生活 (Local Life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value  = 1219
$ws.Range("F10").Value = 869

# --- Sheet: 全部类型 (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value  = 1220
$ws.Range("F8").Value  = 6526
$ws.Range("F10").Value = 752
$ws.Range("F12").Value = 207
$ws.Range("F13").Value = 33
$ws.Range("F14").Value = 762
$ws.Range("F17").Value = 124
$ws.Range("F18").Value = 124
$ws.Range("F21").Value = 83
$ws.Range("F22").Value = 1235
$ws.Range("F23").Value = 17
$ws.Range("F24").Value = 212
$ws.Range("F25").Value = 478
$ws.Range("F28").Value = 352
$ws.Range("F29").Value = 324
$ws.Range("F30").Value = 1436
$ws.Range("F32").Value = 692
$ws.Range("F33").Value = 412
$ws.Range("F34").Value = 416
$ws.Range("F36").Value = 36
$ws.Range("F37").Value = 209
$ws.Range("F44").Value = 133
$ws.Range("F45").Value = 413
$ws.Range("F47").Value = 3679
$ws.Range("F48").Value = 7
$ws.Range("F50").Value = 50
$ws.Range("F51").Value = 672

$wb.Save()
